$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (Topic, Subtopic, Link, File Name)
$ws.Range("A23").Value = "Artificial Intelligence"
$ws.Range("B23").Value = "All the subtopics"
$ws.Range("C23").Value = "https://console.bluemix.net/catalog/?category=ai"
$ws.Range("D23").Value = "AI"

$ws.Range("A25").Value = "Analytics"
$ws.Range("B25").Value = "All the subtopics"
$ws.Range("C25").Value = "https://console.bluemix.net/catalog/?category=analytics"
$ws.Range("D25").Value = "analytics"

$ws.Range("A27").Value = "Databases"
$ws.Range("B27").Value = "All the subtopics"
$ws.Range("C27").Value = "https://console.bluemix.net/catalog/?category=databases"
$ws.Range("D27").Value = "databases"

$ws.Range("A29").Value = "Developer tools"
$ws.Range("B29").Value = "All the subtopics"
$ws.Range("C29").Value = "https://console.bluemix.net/catalog/?category=devops"
$ws.Range("D29").Value = "developer-tools"

$ws.Range("A31").Value = "Integration"
$ws.Range("B31").Value = "All the subtopics"
$ws.Range("C31").Value = "https://console.bluemix.net/catalog/?category=integration"
$ws.Range("D31").Value = "integration"

$ws.Range("A33").Value = "Internet f things"
$ws.Range("B33").Value = "All the subtopics"
$ws.Range("C33").Value = "https://console.bluemix.net/catalog/?category=iot"
$ws.Range("D33").Value = "IOT"

$ws.Range("A35").Value = "Security and Identity"
$ws.Range("B35").Value = "All the subtopics"
$ws.Range("C35").Value = "https://console.bluemix.net/catalog/?category=security"
$ws.Range("D35").Value = "security-identity"

$ws.Range("A37").Value = "Starter kits"
$ws.Range("B37").Value = "All the subtopics"
$ws.Range("C37").Value = "https://console.bluemix.net/catalog/?category=starterkits"
$ws.Range("D37").Value = "starterkits"

$ws.Range("A39").Value = "Web and mobile"
$ws.Range("B39").Value = "All the subtopics"
$ws.Range("C39").Value = "https://console.bluemix.net/catalog/?category=mobile"
$ws.Range("D39").Value = "web-mobile"

$ws.Range("A41").Value = "Web and application"
$ws.Range("B41").Value = "All the subtopics"
$ws.Range("C41").Value = "https://console.bluemix.net/catalog/?category=app_services"
$ws.Range("D41").Value = "web-application"

# Hyperlinks for new rows, then re-apply the same Hyperlink cell style used by C14
$ws.Hyperlinks.Add($ws.Range("C27"), "https://console.bluemix.net/catalog/?category=databases")
$ws.Hyperlinks.Add($ws.Range("C33"), "https://console.bluemix.net/catalog/?category=iot")
$ws.Hyperlinks.Add($ws.Range("C37"), "https://console.bluemix.net/catalog/?category=starterkits")
$ws.Range("C27").Style = $ws.Range("C14").Style
$ws.Range("C33").Style = $ws.Range("C14").Style
$ws.Range("C37").Style = $ws.Range("C14").Style

# Widen column D to fit the new longer text
$ws.Columns.Item(4).ColumnWidth = 16.14

# Update view: zoom level and active selection
$excel.ActiveWindow.Zoom = 75
$ws.Range("C45").Select()
